$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at position 59 (weekly price record), shifting the
# existing rows 59-183 down to 60-184, growing the sheet from 183 to 184 rows.
$ws.Rows(59).Insert(-4121)

# Populate the newly inserted row 59 with the new price record.
$ws.Cells.Item(59, 1).Value = 4
$ws.Cells.Item(59, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(59, 3).Value = "Los Lagos"
$ws.Cells.Item(59, 4).Value = 44571
$ws.Cells.Item(59, 5).Value = 10
$ws.Cells.Item(59, 6).Value = 100112017
$ws.Cells.Item(59, 7).Value = "Apio"
$ws.Cells.Item(59, 8).Value = "Americana (o)"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 25
$ws.Cells.Item(59, 11).Value = 12000
$ws.Cells.Item(59, 12).Value = 12000
$ws.Cells.Item(59, 13).Value = 12000
$ws.Cells.Item(59, 14).Value = "$/docena de matas"
$ws.Cells.Item(59, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(59, 16).Value = 2000
$ws.Cells.Item(59, 17).Value = 6
$ws.Cells.Item(59, 18).Value = "Hortaliza"
